$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D4").Value = "2016-02-16 14:44:57"
$wsZh.Range("G4").Value = "2016-02-16 14:45:55"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D4").Value = "2016-02-16 14:45:14"
$wsDe.Range("G4").Value = "2016-02-16 14:46:24"
